# Update the "想去人数" (interest count) figures in column F for both the
# "展览" and "全部类型" sheets, which hold duplicated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F3"  = 10759
    "F4"  = 600
    "F6"  = 978
    "F8"  = 48
    "F11" = 10568
    "F12" = 4079
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
